# Selenium grid setup: replace sample ID/Pass values with the
# Sumit Deshmukh-prefixed test credentials, and trim the trailing
# blank rows that were left over in the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the header + 2 data rows with the new values
$ws.Range("A1").Value = "SumitDeshmukhIDNO1"
$ws.Range("B1").Value = "SumitDeshmukhPASSNO1"
$ws.Range("A2").Value = "SumitDeshmukhIDNO2"
$ws.Range("B2").Value = "SumitDeshmukhPASSNO2"
$ws.Range("A3").Value = "SumitDeshmukhIDNO3"
$ws.Range("B3").Value = "SumitDeshmukhPASSNO3"

# Remove the now-unused trailing empty rows (4-6) so the used range
# shrinks back down to A1:B3
$ws.Rows("4:6").Delete()

# Select the full data range, matching the new selection state
$ws.Range("A1:B3").Select()
